# Update loading_percent values for the 380 kV case (rows 2-25, cols C-N except B,I,M)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 3).Value = 13.55515853391561
$ws.Cells.Item(2, 4).Value = 6.97220236164387
$ws.Cells.Item(2, 5).Value = 14.04951337151429
$ws.Cells.Item(2, 6).Value = 47.74122866356787
$ws.Cells.Item(2, 7).Value = 59.106604376296
$ws.Cells.Item(2, 8).Value = 21.80101371632457
$ws.Cells.Item(2, 10).Value = 11.46046838249002
$ws.Cells.Item(2, 11).Value = 22.81863010871817
$ws.Cells.Item(2, 12).Value = 9.602360798723558
$ws.Cells.Item(2, 14).Value = 19.73979972027498

# Row 3
$ws.Cells.Item(3, 3).Value = 13.5201211848966
$ws.Cells.Item(3, 4).Value = 6.956957751958257
$ws.Cells.Item(3, 5).Value = 14.05550274076704
$ws.Cells.Item(3, 6).Value = 47.78497829900547
$ws.Cells.Item(3, 7).Value = 59.11881434495974
$ws.Cells.Item(3, 8).Value = 21.85604509383138
$ws.Cells.Item(3, 10).Value = 11.485656650152
$ws.Cells.Item(3, 11).Value = 22.51426913864974
$ws.Cells.Item(3, 12).Value = 9.617832095035554
$ws.Cells.Item(3, 14).Value = 19.81182636252943

# Row 4
$ws.Cells.Item(4, 3).Value = 13.50143597782014
$ws.Cells.Item(4, 4).Value = 6.948259031768182
$ws.Cells.Item(4, 5).Value = 14.06133903944161
$ws.Cells.Item(4, 6).Value = 47.82393102582337
$ws.Cells.Item(4, 7).Value = 59.1439710236677
$ws.Cells.Item(4, 8).Value = 21.89426536496299
$ws.Cells.Item(4, 10).Value = 11.5025733683121
$ws.Cells.Item(4, 11).Value = 22.3300728513259
$ws.Cells.Item(4, 12).Value = 9.628151619452241
$ws.Cells.Item(4, 14).Value = 19.85808573273172

# Row 5
$ws.Cells.Item(5, 3).Value = 13.49453745746508
$ws.Cells.Item(5, 4).Value = 6.944882650258729
$ws.Cells.Item(5, 5).Value = 14.0642607569082
$ws.Cells.Item(5, 6).Value = 47.84283595750721
$ws.Cells.Item(5, 7).Value = 59.15864924047468
$ws.Cells.Item(5, 8).Value = 21.91095154226117
$ws.Cells.Item(5, 10).Value = 11.50983207317018
$ws.Cells.Item(5, 11).Value = 22.25577234720411
$ws.Cells.Item(5, 12).Value = 9.632563564424981
$ws.Cells.Item(5, 14).Value = 19.87745028060218

# Row 6
$ws.Cells.Item(6, 3).Value = 13.4934353248372
$ws.Cells.Item(6, 4).Value = 6.944332235287293
$ws.Cells.Item(6, 5).Value = 14.0647787371288
$ws.Cells.Item(6, 6).Value = 47.84615791028498
$ws.Cells.Item(6, 7).Value = 59.16135343780123
$ws.Cells.Item(6, 8).Value = 21.91378928581467
$ws.Cells.Item(6, 10).Value = 11.51105942482521
$ws.Cells.Item(6, 11).Value = 22.24348326748144
$ws.Cells.Item(6, 12).Value = 9.633308658812984
$ws.Cells.Item(6, 14).Value = 19.88069681580672

# Row 7
$ws.Cells.Item(7, 3).Value = 13.50134003771272
$ws.Cells.Item(7, 4).Value = 6.948212812114732
$ws.Cells.Item(7, 5).Value = 14.06137624203694
$ws.Cells.Item(7, 6).Value = 47.82417372324525
$ws.Cells.Item(7, 7).Value = 59.1441510768347
$ws.Cells.Item(7, 8).Value = 21.89448590596544
$ws.Cells.Item(7, 10).Value = 11.50266978368717
$ws.Cells.Item(7, 11).Value = 22.32906761436751
$ws.Cells.Item(7, 12).Value = 9.628210283168373
$ws.Cells.Item(7, 14).Value = 19.85834480834696

# Row 8
$ws.Cells.Item(8, 3).Value = 13.54249379256525
$ws.Cells.Item(8, 4).Value = 6.966809987444512
$ws.Cells.Item(8, 5).Value = 14.05113080556871
$ws.Cells.Item(8, 6).Value = 47.75379929645138
$ws.Cells.Item(8, 7).Value = 59.10714157759809
$ws.Cells.Item(8, 8).Value = 21.81906712676032
$ws.Cells.Item(8, 10).Value = 11.4688522200073
$ws.Cells.Item(8, 11).Value = 22.71318185261472
$ws.Cells.Item(8, 12).Value = 9.607525355894387
$ws.Cells.Item(8, 14).Value = 19.76421330556849

# Row 9
$ws.Cells.Item(9, 3).Value = 13.64539539956284
$ws.Cells.Item(9, 4).Value = 7.008433976988712
$ws.Cells.Item(9, 5).Value = 14.04814099364365
$ws.Cells.Item(9, 6).Value = 47.71206430949209
$ws.Cells.Item(9, 7).Value = 59.17517968676852
$ws.Cells.Item(9, 8).Value = 21.70645999091447
$ws.Cells.Item(9, 10).Value = 11.41404349545936
$ws.Cells.Item(9, 11).Value = 23.4838769730922
$ws.Cells.Item(9, 12).Value = 9.573449986048155
$ws.Cells.Item(9, 14).Value = 19.59568202394858

# Row 10
$ws.Cells.Item(10, 3).Value = 13.73416899806346
$ws.Cells.Item(10, 4).Value = 7.042027465134167
$ws.Cells.Item(10, 5).Value = 14.05632621578001
$ws.Cells.Item(10, 6).Value = 47.74046867665305
$ws.Cells.Item(10, 7).Value = 59.31136269113988
$ws.Cells.Item(10, 8).Value = 21.64541617051756
$ws.Cells.Item(10, 10).Value = 11.38078313855567
$ws.Cells.Item(10, 11).Value = 24.05559784303366
$ws.Cells.Item(10, 12).Value = 9.552343825957262
$ws.Cells.Item(10, 14).Value = 19.48153390194241

# Row 11
$ws.Cells.Item(11, 3).Value = 13.77732396610976
$ws.Cells.Item(11, 4).Value = 7.057934554938399
$ws.Cells.Item(11, 5).Value = 14.06229171542898
$ws.Cells.Item(11, 6).Value = 47.76626623088655
$ws.Cells.Item(11, 7).Value = 59.39206438176539
$ws.Cells.Item(11, 8).Value = 21.62239011122349
$ws.Cells.Item(11, 10).Value = 11.36717211320776
$ws.Cells.Item(11, 11).Value = 24.31580296121825
$ws.Cells.Item(11, 12).Value = 9.543589705478251
$ws.Cells.Item(11, 14).Value = 19.43167978486472

# Row 12
$ws.Cells.Item(12, 3).Value = 13.7940554818854
$ws.Cells.Item(12, 4).Value = 7.064045396602865
$ws.Cells.Item(12, 5).Value = 14.06487170663258
$ws.Cells.Item(12, 6).Value = 47.77788792652888
$ws.Cells.Item(12, 7).Value = 59.42531712024761
$ws.Cells.Item(12, 8).Value = 21.61435538347439
$ws.Cells.Item(12, 10).Value = 11.36223630639335
$ws.Cells.Item(12, 11).Value = 24.4142660704201
$ws.Cells.Item(12, 12).Value = 9.540396126038834
$ws.Cells.Item(12, 14).Value = 19.41309744974447

# Row 13
$ws.Cells.Item(13, 3).Value = 13.79043486740139
$ws.Cells.Item(13, 4).Value = 7.062725484967059
$ws.Cells.Item(13, 5).Value = 14.06430180772645
$ws.Cells.Item(13, 6).Value = 47.77530258572544
$ws.Cells.Item(13, 7).Value = 59.41803588251108
$ws.Cells.Item(13, 8).Value = 21.61605531073934
$ws.Cells.Item(13, 10).Value = 11.36328960891992
$ws.Cells.Item(13, 11).Value = 24.39306482611249
$ws.Cells.Item(13, 12).Value = 9.541078526988855
$ws.Cells.Item(13, 14).Value = 19.41708633493392

# Row 14
$ws.Cells.Item(14, 3).Value = 13.77869272132698
$ws.Cells.Item(14, 4).Value = 7.058435566469991
$ws.Cells.Item(14, 5).Value = 14.0624975485139
$ws.Cells.Item(14, 6).Value = 47.76718522289891
$ws.Cells.Item(14, 7).Value = 59.39474617114499
$ws.Cells.Item(14, 8).Value = 21.6217153512951
$ws.Cells.Item(14, 10).Value = 11.36676166495676
$ws.Cells.Item(14, 11).Value = 24.3239055039198
$ws.Cells.Item(14, 12).Value = 9.543324536751536
$ws.Cells.Item(14, 14).Value = 19.43014507548226

# Row 15
$ws.Cells.Item(15, 3).Value = 13.77155079083967
$ws.Cells.Item(15, 4).Value = 7.055819137260956
$ws.Cells.Item(15, 5).Value = 14.06143414468653
$ws.Cells.Item(15, 6).Value = 47.76245438811488
$ws.Cells.Item(15, 7).Value = 59.38083105851197
$ws.Cells.Item(15, 8).Value = 21.62527154005345
$ws.Cells.Item(15, 10).Value = 11.36891683834073
$ws.Cells.Item(15, 11).Value = 24.28153158111876
$ws.Cells.Item(15, 12).Value = 9.544716082441633
$ws.Cells.Item(15, 14).Value = 19.43818247347823

# Row 16
$ws.Cells.Item(16, 3).Value = 13.73140374350907
$ws.Cells.Item(16, 4).Value = 7.0410002487939
$ws.Cells.Item(16, 5).Value = 14.05598135444414
$ws.Cells.Item(16, 6).Value = 47.73904198357049
$ws.Cells.Item(16, 7).Value = 59.30646573832639
$ws.Cells.Item(16, 8).Value = 21.64701665452935
$ws.Cells.Item(16, 10).Value = 11.38170321639325
$ws.Cells.Item(16, 11).Value = 24.03858770254107
$ws.Cells.Item(16, 12).Value = 9.55293294090424
$ws.Cells.Item(16, 14).Value = 19.48483354473176

# Row 17
$ws.Cells.Item(17, 3).Value = 13.70747854157206
$ws.Cells.Item(17, 4).Value = 7.032067559335182
$ws.Cells.Item(17, 5).Value = 14.05320952506607
$ws.Cells.Item(17, 6).Value = 47.72797848043589
$ws.Cells.Item(17, 7).Value = 59.26564620597025
$ws.Cells.Item(17, 8).Value = 21.66157336336643
$ws.Cells.Item(17, 10).Value = 11.3899363160251
$ws.Cells.Item(17, 11).Value = 23.88951778886093
$ws.Cells.Item(17, 12).Value = 9.558190410357357
$ws.Cells.Item(17, 14).Value = 19.5139820680656

# Row 18
$ws.Cells.Item(18, 3).Value = 13.69397882579855
$ws.Cells.Item(18, 4).Value = 7.026988732659182
$ws.Cells.Item(18, 5).Value = 14.0518263652452
$ws.Cells.Item(18, 6).Value = 47.72282705827636
$ws.Cells.Item(18, 7).Value = 59.2439333004962
$ws.Cells.Item(18, 8).Value = 21.67039228814916
$ws.Cells.Item(18, 10).Value = 11.39481478390764
$ws.Cells.Item(18, 11).Value = 23.80379382107395
$ws.Cells.Item(18, 12).Value = 9.561294132479571
$ws.Cells.Item(18, 14).Value = 19.53094269447538

# Row 19
$ws.Cells.Item(19, 3).Value = 13.68945321461553
$ws.Cells.Item(19, 4).Value = 7.025279348012263
$ws.Cells.Item(19, 5).Value = 14.05139435185934
$ws.Cells.Item(19, 6).Value = 47.7212909778461
$ws.Cells.Item(19, 7).Value = 59.23688494841065
$ws.Cells.Item(19, 8).Value = 21.67345479058085
$ws.Cells.Item(19, 10).Value = 11.39649111307026
$ws.Cells.Item(19, 11).Value = 23.77477484744058
$ws.Cells.Item(19, 12).Value = 9.562358711524602
$ws.Cells.Item(19, 14).Value = 19.53671884256008

# Row 20
$ws.Cells.Item(20, 3).Value = 13.70999842570671
$ws.Cells.Item(20, 4).Value = 7.033012370582028
$ws.Cells.Item(20, 5).Value = 14.05348275028922
$ws.Cells.Item(20, 6).Value = 47.72903075750833
$ws.Cells.Item(20, 7).Value = 59.26980881751965
$ws.Cells.Item(20, 8).Value = 21.65997756759858
$ws.Cells.Item(20, 10).Value = 11.38904508842791
$ws.Cells.Item(20, 11).Value = 23.90538541823849
$ws.Cells.Item(20, 12).Value = 9.55762249137217
$ws.Cells.Item(20, 14).Value = 19.51085897127676

# Row 21
$ws.Cells.Item(21, 3).Value = 13.78213117235671
$ws.Cells.Item(21, 4).Value = 7.05969327591706
$ws.Cells.Item(21, 5).Value = 14.06301880381132
$ws.Cells.Item(21, 6).Value = 47.76951920598638
$ws.Cells.Item(21, 7).Value = 59.40151388980838
$ws.Cells.Item(21, 8).Value = 21.62003425589126
$ws.Cells.Item(21, 10).Value = 11.36573591146113
$ws.Cells.Item(21, 11).Value = 24.34422189661117
$ws.Cells.Item(21, 12).Value = 9.542661537291096
$ws.Cells.Item(21, 14).Value = 19.42630137760982

# Row 22
$ws.Cells.Item(22, 3).Value = 13.83154084776636
$ws.Cells.Item(22, 4).Value = 7.077637568994523
$ws.Cells.Item(22, 5).Value = 14.07112134206484
$ws.Cells.Item(22, 6).Value = 47.80677918001727
$ws.Cells.Item(22, 7).Value = 59.50328205398561
$ws.Cells.Item(22, 8).Value = 21.59792131527769
$ws.Cells.Item(22, 10).Value = 11.35177493777147
$ws.Cells.Item(22, 11).Value = 24.63057392332007
$ws.Cells.Item(22, 12).Value = 9.533591205121208
$ws.Cells.Item(22, 14).Value = 19.37276461490062

# Row 23
$ws.Cells.Item(23, 3).Value = 13.80496564251095
$ws.Cells.Item(23, 4).Value = 7.068014916779985
$ws.Cells.Item(23, 5).Value = 14.06662625419176
$ws.Cells.Item(23, 6).Value = 47.78590482829672
$ws.Cells.Item(23, 7).Value = 59.44753283341336
$ws.Cells.Item(23, 8).Value = 21.60935726501087
$ws.Cells.Item(23, 10).Value = 11.35910972724434
$ws.Cells.Item(23, 11).Value = 24.4778127184148
$ws.Cells.Item(23, 12).Value = 9.538367608913381
$ws.Cells.Item(23, 14).Value = 19.40118076628754

# Row 24
$ws.Cells.Item(24, 3).Value = 13.7088583905269
$ws.Cells.Item(24, 4).Value = 7.032585044624197
$ws.Cells.Item(24, 5).Value = 14.05335856969799
$ws.Cells.Item(24, 6).Value = 47.72855125677886
$ws.Cells.Item(24, 7).Value = 59.26792143458027
$ws.Cells.Item(24, 8).Value = 21.66069762508216
$ws.Cells.Item(24, 10).Value = 11.38944756029652
$ws.Cells.Item(24, 11).Value = 23.89821172371321
$ws.Cells.Item(24, 12).Value = 9.557878994797482
$ws.Cells.Item(24, 14).Value = 19.51227029163155

# Row 25
$ws.Cells.Item(25, 3).Value = 13.61521586687014
$ws.Cells.Item(25, 4).Value = 6.996634781172009
$ws.Cells.Item(25, 5).Value = 14.04712284457205
$ws.Cells.Item(25, 6).Value = 47.71300332100576
$ws.Cells.Item(25, 7).Value = 59.14165783393744
$ws.Cells.Item(25, 8).Value = 21.73312659351938
$ws.Cells.Item(25, 10).Value = 11.42763952406557
$ws.Cells.Item(25, 11).Value = 23.27404870949315
$ws.Cells.Item(25, 12).Value = 9.58197641619919
$ws.Cells.Item(25, 14).Value = 19.63956719667039
